# Generate Report for Handback
#
# Row 2 on the "zh-cn" and "de-de" sheets (the a3d586ea... file) previously
# shared its Correspond Handoff/Handback datetimes with row 3 (the
# ca6e25c9... file). Re-run the handback report so row 2 gets its own,
# distinct handoff/handback timestamps while row 3 keeps the original ones.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-23 16:54:47"
$ws_zhcn.Range("H2").Value = "2016-03-23 16:55:14"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-23 16:54:51"
$ws_dede.Range("H2").Value = "2016-03-23 16:55:21"
